$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 244-261 currently carry the old "blank filler" style (s=3).
# The target keeps the same style used by the rest of the block (s=4, as seen on A226:B226).
$ws.Range("A226:B226").Copy()
$ws.Range("A244:B261").PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Set-TextValue($cellRef, $text) {
    $ws.Range("ZZ1").Formula = '="' + $text + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

Set-TextValue "A226" "80267046"
$ws.Range("B226").Value = "30259-OSR-I"
$ws.Range("C226").Value = 42000

Set-TextValue "A227" "80267048"
$ws.Range("B227").Value = "23500-GPB-I"
$ws.Range("C227").Value = 2000

Set-TextValue "A228" "80267049"
$ws.Range("B228").Value = "10662-BLB-I"
$ws.Range("C228").Value = 7

Set-TextValue "A229" "80267050"
$ws.Range("B229").Value = "20935-CTY-I"
$ws.Range("C229").Value = 1

Set-TextValue "A230" "80267051"
$ws.Range("B230").Value = "10378-ARI-I"
$ws.Range("C230").Value = 1

Set-TextValue "A231" "80267052"
$ws.Range("B231").Value = "10025-ARI-I"
$ws.Range("C231").Value = 1

Set-TextValue "A232" "80267053"
$ws.Range("B232").Value = "10065-ARI-I"
$ws.Range("C232").Value = 1

Set-TextValue "A233" "80267054"
$ws.Range("B233").Value = "10060-ARI-I"
$ws.Range("C233").Value = 1

Set-TextValue "A234" "80267059"
$ws.Range("B234").Value = "10359-ARI-I"
$ws.Range("C234").Value = 1

Set-TextValue "A235" "80267061"
$ws.Range("B235").Value = "10493-ARI-I"
$ws.Range("C235").Value = 1

Set-TextValue "A236" "80267062"
$ws.Range("B236").Value = "30152-OSR-I"
$ws.Range("C236").Value = 50000

Set-TextValue "A237" "80267062"
$ws.Range("B237").Value = "30173-OSR-L"
$ws.Range("C237").Value = 120000

Set-TextValue "A238" "80267065"
$ws.Range("B238").Value = "10499-ARI-I"
$ws.Range("C238").Value = 1

Set-TextValue "A239" "80267065"
$ws.Range("B239").Value = "10195-ARI-I"
$ws.Range("C239").Value = 2

Set-TextValue "A240" "80267065"
$ws.Range("B240").Value = "10251-ARI-I"
$ws.Range("C240").Value = 1

Set-TextValue "A241" "80267067"
$ws.Range("B241").Value = "11724-DLO-L"
$ws.Range("C241").Value = 9

Set-TextValue "A242" "80267067"
$ws.Range("B242").Value = "15390-DLO-I"
$ws.Range("C242").Value = 1

Set-TextValue "A243" "80267070"
$ws.Range("B243").Value = "10045-ARI-I"
$ws.Range("C243").Value = 1

Set-TextValue "A244" "80267070"
$ws.Range("B244").Value = "10399-ARI-I"
$ws.Range("C244").Value = 1

Set-TextValue "A245" "80267070"
$ws.Range("B245").Value = "10025-ARI-I"
$ws.Range("C245").Value = 1

Set-TextValue "A246" "80267072"
$ws.Range("B246").Value = "10185-ARI-I"
$ws.Range("C246").Value = 1

Set-TextValue "A247" "80267073"
$ws.Range("B247").Value = "10547-ARI-I"
$ws.Range("C247").Value = 1

Set-TextValue "A248" "80267074"
$ws.Range("B248").Value = "10359-ARI-I"
$ws.Range("C248").Value = 2

Set-TextValue "A249" "80267074"
$ws.Range("B249").Value = "10381-ARI-I"
$ws.Range("C249").Value = 2

Set-TextValue "A250" "80267074"
$ws.Range("B250").Value = "10259-ARI-I"
$ws.Range("C250").Value = 2

Set-TextValue "A251" "80267074"
$ws.Range("B251").Value = "10055-ARI-I"
$ws.Range("C251").Value = 2

Set-TextValue "A252" "80267074"
$ws.Range("B252").Value = "10150-ARI-I"
$ws.Range("C252").Value = 1

Set-TextValue "A253" "80267074"
$ws.Range("B253").Value = "10200-ARI-I"
$ws.Range("C253").Value = 1

Set-TextValue "A254" "80267074"
$ws.Range("B254").Value = "10250-ARI-I"
$ws.Range("C254").Value = 1

Set-TextValue "A255" "84004838"
$ws.Range("B255").Value = "20953-CTY-I"
$ws.Range("C255").Value = 3

Set-TextValue "A256" "84004839"
$ws.Range("B256").Value = "10257-ARI-I"
$ws.Range("C256").Value = 1

Set-TextValue "A257" "84004840"
$ws.Range("B257").Value = "10355-ARI-I"
$ws.Range("C257").Value = 1

Set-TextValue "A258" "84004841"
$ws.Range("B258").Value = "10255-ARI-I"
$ws.Range("C258").Value = 1

Set-TextValue "A259" "84004842"
$ws.Range("B259").Value = "10587-ARI-I"
$ws.Range("C259").Value = 1

Set-TextValue "A260" "84004843"
$ws.Range("B260").Value = "10000-LDG-I"
$ws.Range("C260").Value = 1

Set-TextValue "A261" "84004844"
$ws.Range("B261").Value = "10040-ARI-I"
$ws.Range("C261").Value = 1

$ws.Range("ZZ1").Clear()

$ws.Range("A15:C261").Select()